# Applies the commit "filter date et restriction sur l'ajout":
#  - replaces the content of rows 2-3 with the new "equine encephalomyelitis" entries
#  - appends 3 new rows (4-6) of scraped disease-alert data
#  - the Date column (D) is stored as plain text (e.g. "2023-12-04"),
#    not as a parsed Excel date serial, so it survives any downstream date filter
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ A = 'equine encephalomyelitis'; B = 'Equine encephalomyelitis, Argentina declares sanitary emergency; Uruguay reports first case Equine encephalomyelitis, Argentina declares sanitary emergency; Uruguay reports first case Equine encephalomyelitis is highly contagious and can affect humans Argentina has finally declared a sanitary emergency in all of its territory, following an increase in equine encephalomyelitis, EE, and anticipates immediate, extraordinary, exceptional measures to contain further outbreaks of the virus disease. The Argentine government points out that the East, West and Venezuelan variants of EE, are viral infections, transmitted by mosquitoes, which can cause very serious encephalitis in horses and eventually humans.'; C = 'https://en.mercopress.com/2023/12/04/equine-encephalomyelitis-argentina-declares-sanitary-emergency-uruguay-reports-first-case'; D = '2023-12-04'; E = 'Paysandu, Argentina' },
  @{ A = 'equine encephalomyelitis'; B = 'Equine encephalomyelitis, Argentina declares sanitary emergency; Uruguay reports first case Equine encephalomyelitis, Argentina declares sanitary emergency; Uruguay reports first case Equine encephalomyelitis is highly contagious and can affect humans Argentina has finally declared a sanitary emergency in all of its territory, following an increase in equine encephalomyelitis, EE, and anticipates immediate, extraordinary, exceptional measures to contain further outbreaks of the virus disease. The Argentine government points out that the East, West and Venezuelan variants of EE, are viral infections, transmitted by mosquitoes, which can cause very serious encephalitis in horses and eventually humans.'; C = 'https://en.mercopress.com/2023/12/04/equine-encephalomyelitis-argentina-declares-sanitary-emergency-uruguay-reports-first-case'; D = '2023-12-04'; E = 'Paysandu, Argentina' },
  @{ A = 'bluetongue'; B = 'But the disease has resurfaced this year and animal health experts reported that a new strain, named BTV-3, had been confirmed on more than 700 Dutch farms by early October - with a new case also reported in Belgium this week. “The existing BTV-8 serotype vaccine will not offer cross-protection against this new BTV-3 strain, making any likely outbreak difficult to control. "Hence why it’s so important that we follow the advice to take action and prioritise good biosecurity measures while remaining extremely vigilant to the disease at this stage." Dr Henry added: "It remains extremely difficult to protect against midges and a vector-borne disease.'; C = 'https://lc.cx/nLtrtC'; D = '2023-10-11'; E = 'Belgium' },
  @{ A = 'brucellosis'; B = 'Brucellosis, a zoonotic disease caused by the bacterial genus Brucella, has been confirmed in two natives of Vembayam in Thiruvananthapuram. A statement issued by the Animal Husbandry department said that it was difficult to recognise the disease in animals, as it did not produce any overt symptoms in animals. The bacteria are transmitted from animals to humans by ingestion through infected food products, direct contact with an infected animal, or through the inhalation of aerosols. Minister for Animal Husbandry J. Chinchurani, said that apart from giving awareness classes to dairy farmers, the department would test milk samples from milk societies also.'; C = 'https://lc.cx/4H9k-A'; D = '2023-10-09'; E = 'Thiruvananthapuram, Kerala, Kollam' },
  @{ A = 'west nile virus, eastern equine encephalitis'; B = 'The first EEE positive mosquito pool was detected in Gloucester County this year (https://www.nj.gov/health/cd/statistics/arboviral-stats/). “Vaccinated animals are much less likely to contract deadly diseases such as EEE and West Nile Virus.” For more information about EEE in horses, visit the New Jersey Department of Agriculture website at: http://www.nj.gov/agriculture/divisions/ah/diseases/diseaseworksheets.html EEE and West Nile virus, like other viral diseases affecting a horse’s neurological system, must be reported to the state veterinarian at 609-671-6400 within 48 hours of diagnosis.'; C = 'https://lc.cx/RqYUj8'; D = '2023-10-25'; E = 'New Jersey, Gloucester County' }
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row.A
  $ws.Cells.Item($r, 2).Value = $row.B
  $ws.Cells.Item($r, 3).Value = $row.C

  # Write the date as literal text: stamp the cell Text ("@") before
  # assigning so Excel does not auto-convert the "yyyy-mm-dd" string into
  # a date serial, then drop back to the default "Normal" style so the
  # cell carries no leftover number-format override (plain inline string).
  $d = $ws.Cells.Item($r, 4)
  $d.Style = "Normal"
  $d.NumberFormat = "@"
  $d.Value = $row.D
  $d.Style = "Normal"

  $ws.Cells.Item($r, 5).Value = $row.E

  $r = $r + 1
}
